# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2410"
#   "<header>_new" -> "<header>_FV2504"
# then format the header row + data as an Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21
$lastRow = 58

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = $cell.Value()
    if ($old -like "*_old") {
        $cell.Value = $old -replace "_old$", "_FV2410"
    } elseif ($old -like "*_new") {
        $cell.Value = $old -replace "_new$", "_FV2504"
    }
}

# Turn the range into a proper Excel Table (Table1) spanning the used range.
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1, top-left cell of the scrolling
# pane is A2) and keep the sheet's own selection anchored there.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
